$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 26 de Julio de 2020 a las 01:41"

# Swap country names to reflect reordering (Guayana Francesa before Haiti,
# Luxemburgo before Mauritania)
$ws.Cells.Item(90, 1).Value = "Guayana Francesa"
$ws.Cells.Item(91, 1).Value = "Haiti"
$ws.Cells.Item(95, 1).Value = "Luxemburgo"
$ws.Cells.Item(96, 1).Value = "Mauritania"

# Update statistic values per country row
# Row 4
$ws.Cells.Item(4, 2).Value = 4312284
$ws.Cells.Item(4, 3).Value = 63957
$ws.Cells.Item(4, 4).Value = 2055731
$ws.Cells.Item(4, 5).Value = 2107213
$ws.Cells.Item(4, 7).Value = 850
$ws.Cells.Item(4, 8).Value = 149340

# Row 5
$ws.Cells.Item(5, 2).Value = 2396434
$ws.Cells.Item(5, 3).Value = 48234
$ws.Cells.Item(5, 4).Value = 1617480
$ws.Cells.Item(5, 5).Value = 692458
$ws.Cells.Item(5, 7).Value = 1111
$ws.Cells.Item(5, 8).Value = 86496

# Row 15
$ws.Cells.Item(15, 2).Value = 273112
$ws.Cells.Item(15, 3).Value = 2712
$ws.Cells.Item(15, 4).Value = 237434
$ws.Cells.Item(15, 5).Value = 29856
$ws.Cells.Item(15, 7).Value = 59
$ws.Cells.Item(15, 8).Value = 5822

# Row 18
$ws.Cells.Item(18, 2).Value = 240795
$ws.Cells.Item(18, 3).Value = 7254
$ws.Cells.Item(18, 4).Value = 119667
$ws.Cells.Item(18, 5).Value = 112859
$ws.Cells.Item(18, 7).Value = 294
$ws.Cells.Item(18, 8).Value = 8269

# Row 23
$ws.Cells.Item(23, 2).Value = 158334
$ws.Cells.Item(23, 3).Value = 4814
$ws.Cells.Item(23, 5).Value = 87419
$ws.Cells.Item(23, 7).Value = 86
$ws.Cells.Item(23, 8).Value = 2893

# Row 24
$ws.Cells.Item(24, 2).Value = 113556
$ws.Cells.Item(24, 3).Value = 350
$ws.Cells.Item(24, 4).Value = 99125
$ws.Cells.Item(24, 5).Value = 5546

# Row 50
$ws.Cells.Item(50, 2).Value = 39977
$ws.Cells.Item(50, 3).Value = 438
$ws.Cells.Item(50, 4).Value = 16948
$ws.Cells.Item(50, 5).Value = 22173
$ws.Cells.Item(50, 7).Value = 11
$ws.Cells.Item(50, 8).Value = 856

# Row 57
$ws.Cells.Item(57, 2).Value = 31851
$ws.Cells.Item(57, 3).Value = 794
$ws.Cells.Item(57, 4).Value = 28438
$ws.Cells.Item(57, 5).Value = 3252

# Row 59
$ws.Cells.Item(59, 2).Value = 28786
$ws.Cells.Item(59, 3).Value = 830
$ws.Cells.Item(59, 4).Value = 21567
$ws.Cells.Item(59, 5).Value = 6226
$ws.Cells.Item(59, 7).Value = 1
$ws.Cells.Item(59, 8).Value = 993

# Row 71
$ws.Cells.Item(71, 2).Value = 15212
$ws.Cells.Item(71, 3).Value = 131
$ws.Cells.Item(71, 5).Value = 5253

# Row 85
$ws.Cells.Item(85, 2).Value = 9111
$ws.Cells.Item(85, 3).Value = 19
$ws.Cells.Item(85, 5).Value = 182

# Row 90
$ws.Cells.Item(90, 2).Value = 7332
$ws.Cells.Item(90, 3).Value = 81
$ws.Cells.Item(90, 4).Value = 5767
$ws.Cells.Item(90, 5).Value = 1524
$ws.Cells.Item(90, 7).Value = 0
$ws.Cells.Item(90, 8).Value = 41

# Row 91
$ws.Cells.Item(91, 2).Value = 7260
$ws.Cells.Item(91, 3).Value = 63
$ws.Cells.Item(91, 4).Value = 4236
$ws.Cells.Item(91, 5).Value = 2868
$ws.Cells.Item(91, 7).Value = 2
$ws.Cells.Item(91, 8).Value = 156

# Row 94
$ws.Cells.Item(94, 2).Value = 6927
$ws.Cells.Item(94, 3).Value = 60
$ws.Cells.Item(94, 4).Value = 6098
$ws.Cells.Item(94, 5).Value = 787

# Row 95
$ws.Cells.Item(95, 2).Value = 6189
$ws.Cells.Item(95, 3).Value = 133
$ws.Cells.Item(95, 4).Value = 4647
$ws.Cells.Item(95, 5).Value = 1430
$ws.Cells.Item(95, 8).Value = 112

# Row 96
$ws.Cells.Item(96, 2).Value = 6151
$ws.Cells.Item(96, 3).Value = 35
$ws.Cells.Item(96, 4).Value = 4299
$ws.Cells.Item(96, 5).Value = 1696
$ws.Cells.Item(96, 8).Value = 156

# Row 100
$ws.Cells.Item(100, 2).Value = 4598
$ws.Cells.Item(100, 3).Value = 5
$ws.Cells.Item(100, 4).Value = 1506
$ws.Cells.Item(100, 5).Value = 3033

# Row 103
$ws.Cells.Item(103, 2).Value = 4328
$ws.Cells.Item(103, 3).Value = 104
$ws.Cells.Item(103, 4).Value = 2679
$ws.Cells.Item(103, 5).Value = 1609
$ws.Cells.Item(103, 7).Value = 2
$ws.Cells.Item(103, 8).Value = 40

# Row 115
$ws.Cells.Item(115, 2).Value = 2747
$ws.Cells.Item(115, 3).Value = 82
$ws.Cells.Item(115, 4).Value = 664
$ws.Cells.Item(115, 5).Value = 2040

# Row 138
$ws.Cells.Item(138, 2).Value = 1381
$ws.Cells.Item(138, 3).Value = 76
$ws.Cells.Item(138, 4).Value = 853
$ws.Cells.Item(138, 5).Value = 505

# Row 140
$ws.Cells.Item(140, 2).Value = 1174
$ws.Cells.Item(140, 3).Value = 8
$ws.Cells.Item(140, 4).Value = 947
$ws.Cells.Item(140, 5).Value = 193

# Row 143
$ws.Cells.Item(143, 4).Value = 1025
$ws.Cells.Item(143, 5).Value = 30

# Row 158
$ws.Cells.Item(158, 2).Value = 657
$ws.Cells.Item(158, 3).Value = 3
$ws.Cells.Item(158, 5).Value = 96

# Row 183
$ws.Cells.Item(183, 2).Value = 119
$ws.Cells.Item(183, 3).Value = 1
$ws.Cells.Item(183, 4).Value = 102
